$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "59.494.92"
$ws.Range("E2").Value = "  +1.99%  "

$ws.Range("D3").Value = "2.610.56"
$ws.Range("E3").Value = "  +1.54%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = "  -0.18%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "537.62"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +4.01%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "141.54"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  +2.84%  "

$ws.Range("E7").Value = "  +0.11%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.570"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  +1.24%  "

$ws.Range("D9").Value = "2.624.33"
$ws.Range("E9").Value = "  +1.51%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "6.49"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  +0.05%  "

$ws.Range("E11").Value = "  +4.27%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.338"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  +3.51%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.135"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  +2.05%  "

$ws.Range("D14").Value = "3.068.98"
$ws.Range("E14").Value = "  +1.35%  "

$ws.Range("D15").Value = "59.420.16"
$ws.Range("E15").Value = "  +1.93%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "20.63"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  +2.35%  "

$ws.Range("D17").Value = "2.626.39"
$ws.Range("E17").Value = "  +2.10%  "

$ws.Range("E18").Value = "  +2.67%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "346.12"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  +3.31%  "

$ws.Range("E20").Value = "  +1.96%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "10.16"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  +1.25%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.39"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  +0.47%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.999"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  +0.08%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "67.14"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  +2.17%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.168"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  +1.56%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.409"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  +2.52%  "

$ws.Range("E27").Value = "  +0.25%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "7.23"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  +4.20%  "

$ws.Range("D29").Value = "0.0₃0752"
$ws.Range("E29").Value = "  +7.10%  "

$ws.Range("E30").Value = "  +0.06%  "

$ws.Range("E31").Value = "  +6.10%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "5.88"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  +1.38%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "18.94"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  +1.89%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "149.22"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  +0.55%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "4.02"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  +3.57%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.13"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  +1.80%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "36.94"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  +2.03%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.846"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  +3.12%  "

$ws.Range("E39").Value = "  +3.82%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.843"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  +3.00%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "3.56"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  +2.50%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "278.03"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  +3.36%  "

$ws.Range("E43").Value = "  +0.19%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.602"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  +2.67%  "

$ws.Range("B45").Value = "Stellar"
$ws.Range("C45").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0962"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  +2.36%  "

$ws.Range("B46").Value = "WhiteBITCoin"
$ws.Range("C46").Value = "https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "10.73"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  +0.04%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0525"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  +2.35%  "

$ws.Range("B48").Value = "Maker"
$ws.Range("C48").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D48").Value = "1.956.55"
$ws.Range("E48").Value = "  -0.66%  "

$ws.Range("B49").Value = "VeChain"
$ws.Range("C49").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0224"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  +2.97%  "

$ws.Range("B50").Value = "InjectiveProtocol"
$ws.Range("C50").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "18.43"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  +5.65%  "

$ws.Range("B51").Value = "RenderToken"
$ws.Range("C51").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "4.54"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  +3.17%  "
